$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column R (18) into new column S (19) for rows 3-34,
# mirroring the style Excel applies when a user extends a table by one
# column (drag-fill / copy format right).
$ws.Range("R3:R34").Copy() | Out-Null
$ws.Range("S3:S34").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New data for the 2022 column (S)
$values = @{
    4  = 2022
    5  = 135
    6  = 99
    7  = 36
    8  = 97
    9  = 80
    10 = 17
    11 = 17
    12 = 11
    13 = 6
    14 = 5
    15 = 3
    16 = 2
    17 = "-"
    18 = "-"
    19 = "-"
    20 = 6
    21 = 1
    22 = 5
    23 = "-"
    24 = "-"
    25 = "-"
    26 = 10
    27 = 4
    28 = 6
    29 = "-"
    30 = "-"
    31 = "-"
    32 = "-"
    33 = "-"
    34 = "-"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 19).Value = $values[$row]
}

# Update the active selection to reflect where the user ended up (next to
# the newly filled column), matching the saved sheetView state.
$ws.Range("T4").Select() | Out-Null
